$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Cells.Item(15, 8).Value = 1491.5542
$ws.Cells.Item(15, 9).Value = 1491.5542
$ws.Cells.Item(15, 11).Value = 4474.6626
$ws.Cells.Item(15, 13).Value = -4305.6626

# Row 34
$ws.Cells.Item(34, 8).Value = 1065.4
$ws.Cells.Item(34, 9).Value = 1065.4
$ws.Cells.Item(34, 11).Value = 1065.4
$ws.Cells.Item(34, 13).Value = -862.4000000000001

# Row 36
$ws.Cells.Item(36, 8).Value = 1065.4
$ws.Cells.Item(36, 9).Value = 1065.4
$ws.Cells.Item(36, 11).Value = 1065.4
$ws.Cells.Item(36, 13).Value = -350.4000000000001

# Row 46
$ws.Cells.Item(46, 8).Value = 1159.9
$ws.Cells.Item(46, 10).Value = 1177.6666
$ws.Cells.Item(46, 12).Value = 3532.9998
$ws.Cells.Item(46, 14).Value = -3770.9998

# Row 60
$ws.Cells.Item(60, 8).Value = 1159.9
$ws.Cells.Item(60, 10).Value = 1177.6666
$ws.Cells.Item(60, 12).Value = 3532.9998
$ws.Cells.Item(60, 14).Value = -4500.9998

# Row 137
$ws.Cells.Item(137, 8).Value = 2392.35
$ws.Cells.Item(137, 9).Value = 0
$ws.Cells.Item(137, 10).Value = 2392.35
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 7177.049999999999
$ws.Cells.Item(137, 14).Value = -12277.05
$ws.Cells.Item(137, 13).Value = $null

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Cells.Item(32, 8).Value = 5958416
$ws.Cells.Item(32, 9).Value = 5350.696
$ws.Cells.Item(32, 10).Value = 33342518
$ws.Cells.Item(32, 11).Value = 5350.696
$ws.Cells.Item(32, 12).Value = 33342518
$ws.Cells.Item(32, 13).Value = -5063.696
$ws.Cells.Item(32, 14).Value = -33343092

# Row 110
$ws.Cells.Item(110, 8).Value = 9800
$ws.Cells.Item(110, 9).Value = 0
$ws.Cells.Item(110, 10).Value = 9800
$ws.Cells.Item(110, 11).Value = 0
$ws.Cells.Item(110, 12).Value = 9800
$ws.Cells.Item(110, 14).Value = -13890
$ws.Cells.Item(110, 13).Value = $null

# Row 122
$ws.Cells.Item(122, 8).Value = 1746.5714
$ws.Cells.Item(122, 9).Value = 1746.5714
$ws.Cells.Item(122, 10).Value = 0
$ws.Cells.Item(122, 11).Value = 5239.7142
$ws.Cells.Item(122, 12).Value = 0
$ws.Cells.Item(122, 13).Value = -2789.7142
$ws.Cells.Item(122, 14).Value = $null

$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Cells.Item(94, 8).Value = 400.8421
$ws.Cells.Item(94, 9).Value = 403.88235
$ws.Cells.Item(94, 10).Value = 375
$ws.Cells.Item(94, 11).Value = 403.88235
$ws.Cells.Item(94, 12).Value = 375
$ws.Cells.Item(94, 13).Value = 47.11765000000003
$ws.Cells.Item(94, 14).Value = -1277

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Cells.Item(31, 8).Value = 1160.8684
$ws.Cells.Item(31, 9).Value = 718.625
$ws.Cells.Item(31, 10).Value = 1482.5
$ws.Cells.Item(31, 11).Value = 718.625
$ws.Cells.Item(31, 12).Value = 1482.5
$ws.Cells.Item(31, 14).Value = -2072.5
$ws.Cells.Item(31, 13).Value = -423.625

# Row 34
$ws.Cells.Item(34, 8).Value = 1160.8684
$ws.Cells.Item(34, 9).Value = 718.625
$ws.Cells.Item(34, 10).Value = 1482.5
$ws.Cells.Item(34, 11).Value = 718.625
$ws.Cells.Item(34, 12).Value = 1482.5
$ws.Cells.Item(34, 14).Value = -1886.5
$ws.Cells.Item(34, 13).Value = -516.625

# Row 70
$ws.Cells.Item(70, 8).Value = 40998
$ws.Cells.Item(70, 10).Value = 40998
$ws.Cells.Item(70, 12).Value = 40998
$ws.Cells.Item(70, 14).Value = -41628

# Row 73
$ws.Cells.Item(73, 8).Value = 40998
$ws.Cells.Item(73, 10).Value = 40998
$ws.Cells.Item(73, 12).Value = 40998
$ws.Cells.Item(73, 14).Value = -43182

# Row 99
$ws.Cells.Item(99, 8).Value = 34488612
$ws.Cells.Item(99, 9).Value = 83343320
$ws.Cells.Item(99, 10).Value = 2933.4119
$ws.Cells.Item(99, 11).Value = 83343320
$ws.Cells.Item(99, 12).Value = 2933.4119
$ws.Cells.Item(99, 13).Value = -83341822
$ws.Cells.Item(99, 14).Value = -5929.4119

# Row 126
$ws.Cells.Item(126, 8).Value = 34488612
$ws.Cells.Item(126, 9).Value = 83343320
$ws.Cells.Item(126, 10).Value = 2933.4119
$ws.Cells.Item(126, 11).Value = 250029960
$ws.Cells.Item(126, 12).Value = 8800.235700000001
$ws.Cells.Item(126, 13).Value = -250027490
$ws.Cells.Item(126, 14).Value = -13740.2357

$ws = $wb.Worksheets.Item("CUL")
# Row 3
$ws.Cells.Item(3, 8).Value = 5024.2856
$ws.Cells.Item(3, 9).Value = 3855
$ws.Cells.Item(3, 10).Value = 6583.3335
$ws.Cells.Item(3, 11).Value = 11565
$ws.Cells.Item(3, 12).Value = 19750.0005
$ws.Cells.Item(3, 13).Value = -11453
$ws.Cells.Item(3, 14).Value = -19974.0005

# Row 5
$ws.Cells.Item(5, 8).Value = 24825344
$ws.Cells.Item(5, 9).Value = 43210180
$ws.Cells.Item(5, 10).Value = 5816.75
$ws.Cells.Item(5, 11).Value = 129630540
$ws.Cells.Item(5, 12).Value = 17450.25
$ws.Cells.Item(5, 13).Value = -129630428
$ws.Cells.Item(5, 14).Value = -17674.25

# Row 37
$ws.Cells.Item(37, 8).Value = 722501.75
$ws.Cells.Item(37, 10).Value = 722501.75
$ws.Cells.Item(37, 12).Value = 2167505.25
$ws.Cells.Item(37, 14).Value = -2167729.25

# Row 68
$ws.Cells.Item(68, 8).Value = 3226.4285
$ws.Cells.Item(68, 9).Value = 672.9032
$ws.Cells.Item(68, 10).Value = 10422.728
$ws.Cells.Item(68, 11).Value = 2018.7096
$ws.Cells.Item(68, 12).Value = 31268.184
$ws.Cells.Item(68, 13).Value = -1207.7096
$ws.Cells.Item(68, 14).Value = -32890.18399999999

# Row 71
$ws.Cells.Item(71, 8).Value = 3226.4285
$ws.Cells.Item(71, 9).Value = 672.9032
$ws.Cells.Item(71, 10).Value = 10422.728
$ws.Cells.Item(71, 11).Value = 6056.1288
$ws.Cells.Item(71, 12).Value = 93804.552
$ws.Cells.Item(71, 13).Value = -2000.1288
$ws.Cells.Item(71, 14).Value = -101916.552

# Row 87
$ws.Cells.Item(87, 8).Value = 1086.8572
$ws.Cells.Item(87, 9).Value = 1086.8572
$ws.Cells.Item(87, 11).Value = 3260.5716
$ws.Cells.Item(87, 13).Value = -2012.5716

# Row 90
$ws.Cells.Item(90, 8).Value = 1086.8572
$ws.Cells.Item(90, 9).Value = 1086.8572
$ws.Cells.Item(90, 11).Value = 9781.7148
$ws.Cells.Item(90, 13).Value = -3541.7148

# Row 123
$ws.Cells.Item(123, 8).Value = 1907.5
$ws.Cells.Item(123, 9).Value = 876.6667
$ws.Cells.Item(123, 10).Value = 5000
$ws.Cells.Item(123, 11).Value = 2630.0001
$ws.Cells.Item(123, 12).Value = 15000
$ws.Cells.Item(123, 13).Value = -180.0001000000002
$ws.Cells.Item(123, 14).Value = -19900

# Row 131
$ws.Cells.Item(131, 8).Value = 963.45
$ws.Cells.Item(131, 9).Value = 627.6667
$ws.Cells.Item(131, 10).Value = 996.65936
$ws.Cells.Item(131, 11).Value = 1883.0001
$ws.Cells.Item(131, 12).Value = 2989.97808
$ws.Cells.Item(131, 13).Value = 3156.9999
$ws.Cells.Item(131, 14).Value = -13069.97808

# Row 135
$ws.Cells.Item(135, 8).Value = 24825344
$ws.Cells.Item(135, 9).Value = 43210180
$ws.Cells.Item(135, 10).Value = 5816.75
$ws.Cells.Item(135, 11).Value = 388891620
$ws.Cells.Item(135, 12).Value = 52350.75
$ws.Cells.Item(135, 13).Value = -388889085
$ws.Cells.Item(135, 14).Value = -57420.75

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Cells.Item(15, 8).Value = 39400
$ws.Cells.Item(15, 10).Value = 39400
$ws.Cells.Item(15, 12).Value = 39400
$ws.Cells.Item(15, 14).Value = -39976

# Row 81
$ws.Cells.Item(81, 8).Value = 39400
$ws.Cells.Item(81, 10).Value = 39400
$ws.Cells.Item(81, 12).Value = 39400
$ws.Cells.Item(81, 14).Value = -41396

# Row 84
$ws.Cells.Item(84, 8).Value = 39400
$ws.Cells.Item(84, 10).Value = 39400
$ws.Cells.Item(84, 12).Value = 118200
$ws.Cells.Item(84, 14).Value = -128184
